$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 54 (hunk 0)
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents() | Out-Null
# row 107 (hunk 1)
$ws.Range("H107").Value = 1240.8889
$ws.Range("I107").Value = 1108.875
$ws.Range("K107").Value = 1108.875
$ws.Range("M107").Value = 811.125
# row 113 (hunk 2)
$ws.Range("H113").Value = 2687.5
$ws.Range("I113").Value = 2574.4
$ws.Range("K113").Value = 2574.4
$ws.Range("M113").Value = 679.5999999999999
# row 138 (hunk 3)
$ws.Range("H138").Value = 1966.421
$ws.Range("J138").Value = 2770.4
$ws.Range("L138").Value = 8311.200000000001
$ws.Range("N138").Value = -18591.2
# row 141 (hunk 4)
$ws.Range("H141").Value = 90009.45
$ws.Range("I141").Value = 98730.39999999999
$ws.Range("K141").Value = 296191.2
$ws.Range("M141").Value = -291011.2

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 45 (hunk 5)
$ws.Range("H45").Value = 6250
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents() | Out-Null
# row 122 (hunk 6)
$ws.Range("H122").Value = 2750
$ws.Range("I122").Value = 2750
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8250
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5800
$ws.Range("N122").ClearContents() | Out-Null

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 20 (hunk 7)
$ws.Range("H20").Value = 9003.727999999999
$ws.Range("I20").Value = 11294
$ws.Range("J20").Value = 2896.3333
$ws.Range("K20").Value = 11294
$ws.Range("L20").Value = 2896.3333
$ws.Range("M20").Value = -11047
$ws.Range("N20").Value = -3390.3333
# row 35 (hunk 8)
$ws.Range("H35").Value = 16378
$ws.Range("I35").Value = 3567
$ws.Range("J35").Value = 42000
$ws.Range("K35").Value = 3567
$ws.Range("L35").Value = 42000
$ws.Range("M35").Value = -3257
$ws.Range("N35").Value = -42620
# row 94 (hunk 9)
$ws.Range("H94").Value = 12549.9
$ws.Range("I94").Value = 2728.8572
$ws.Range("K94").Value = 2728.8572
$ws.Range("M94").Value = -2277.8572
# row 107 (hunk 10)
$ws.Range("H107").Value = 1234.9656
$ws.Range("I107").Value = 1105.9565
$ws.Range("J107").Value = 1729.5
$ws.Range("K107").Value = 1105.9565
$ws.Range("L107").Value = 1729.5
$ws.Range("M107").Value = 814.0435
$ws.Range("N107").Value = -5569.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 21 (hunk 11)
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents() | Out-Null
# row 26 (hunk 12)
$ws.Range("H26").Value = 10019.5
$ws.Range("J26").Value = 10019
$ws.Range("L26").Value = 10019
$ws.Range("N26").Value = -10593
# row 31 (hunk 13)
$ws.Range("H31").Value = 2111.6155
$ws.Range("I31").Value = 2055.0833
$ws.Range("K31").Value = 2055.0833
$ws.Range("M31").Value = -1760.0833
# row 33 (hunk 14)
$ws.Range("H33").Value = 12874.375
$ws.Range("I33").Value = 4482.75
$ws.Range("K33").Value = 4482.75
$ws.Range("M33").Value = -4103.75
# row 34 (hunk 15)
$ws.Range("H34").Value = 2111.6155
$ws.Range("I34").Value = 2055.0833
$ws.Range("K34").Value = 2055.0833
$ws.Range("M34").Value = -1853.0833
# row 35 (hunk 16)
$ws.Range("H35").Value = 8129.25
$ws.Range("I35").Value = 829.6667
$ws.Range("J35").Value = 30028
$ws.Range("K35").Value = 829.6667
$ws.Range("L35").Value = 30028
$ws.Range("M35").Value = -535.6667
$ws.Range("N35").Value = -30616
# row 36 (hunk 17)
$ws.Range("H36").Value = 15219.8
$ws.Range("I36").Value = 9348.666999999999
$ws.Range("J36").Value = 24026.5
$ws.Range("K36").Value = 9348.666999999999
$ws.Range("L36").Value = 24026.5
$ws.Range("M36").Value = -8960.666999999999
$ws.Range("N36").Value = -24802.5
# row 40 (hunk 18)
$ws.Range("H40").Value = 15219.8
$ws.Range("I40").Value = 9348.666999999999
$ws.Range("J40").Value = 24026.5
$ws.Range("K40").Value = 9348.666999999999
$ws.Range("L40").Value = 24026.5
$ws.Range("M40").Value = -9188.666999999999
$ws.Range("N40").Value = -24346.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 2 (hunk 19)
$ws.Range("H2").Value = 36.15
$ws.Range("I2").Value = 26.285715
$ws.Range("J2").Value = 41.46154
$ws.Range("K2").Value = 157.71429
$ws.Range("L2").Value = 248.76924
$ws.Range("M2").Value = -44.71429000000001
$ws.Range("N2").Value = -474.76924
# row 11 (hunk 20)
$ws.Range("H11").Value = 219.25
$ws.Range("I11").Value = 100.833336
$ws.Range("J11").Value = 574.5
$ws.Range("K11").Value = 302.500008
$ws.Range("L11").Value = 1723.5
$ws.Range("M11").Value = -162.500008
$ws.Range("N11").Value = -2003.5
# row 26 (hunk 21)
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 50
$ws.Range("K26").Value = 150
$ws.Range("M26").Value = 138
# row 68 (hunk 22)
$ws.Range("H68").Value = 1740.3334
$ws.Range("J68").Value = 1499.5
$ws.Range("L68").Value = 4498.5
$ws.Range("N68").Value = -6120.5
# row 71 (hunk 23)
$ws.Range("H71").Value = 1740.3334
$ws.Range("J71").Value = 1499.5
$ws.Range("L71").Value = 13495.5
$ws.Range("N71").Value = -21607.5
# row 82 (hunk 24)
$ws.Range("H82").Value = 9500
$ws.Range("J82").Value = 9500
$ws.Range("L82").Value = 28500
$ws.Range("N82").Value = -29312
# row 85 (hunk 25)
$ws.Range("H85").Value = 9500
$ws.Range("J85").Value = 9500
$ws.Range("L85").Value = 28500
$ws.Range("N85").Value = -31308

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 126 (hunk 26)
$ws.Range("H126").Value = 13088.556
$ws.Range("J126").Value = 15760
$ws.Range("L126").Value = 47280
$ws.Range("N126").Value = -52220

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 16 (hunk 27)
$ws.Range("H16").Value = 1113.4615
$ws.Range("I16").Value = 1113.4615
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1113.4615
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -943.4614999999999
$ws.Range("N16").ClearContents() | Out-Null
# row 22 (hunk 28)
$ws.Range("H22").Value = 1639.8572
$ws.Range("I22").Value = 2096.2
$ws.Range("K22").Value = 2096.2
$ws.Range("M22").Value = -1801.2
# row 27 (hunk 29)
$ws.Range("H27").Value = 1639.8572
$ws.Range("I27").Value = 2096.2
$ws.Range("K27").Value = 2096.2
$ws.Range("M27").Value = -1989.2
# row 46 (hunk 30)
$ws.Range("H46").Value = 3097.7576
$ws.Range("I46").Value = 700
$ws.Range("J46").Value = 3525.9285
$ws.Range("K46").Value = 700
$ws.Range("L46").Value = 3525.9285
$ws.Range("M46").Value = -512
$ws.Range("N46").Value = -3901.9285
# row 82 (hunk 31)
$ws.Range("H82").Value = 2404.4092
$ws.Range("I82").Value = 2480.8096
$ws.Range("J82").Value = 800
$ws.Range("K82").Value = 2480.8096
$ws.Range("L82").Value = 800
$ws.Range("M82").Value = -2119.8096
$ws.Range("N82").Value = -1522
# row 85 (hunk 32)
$ws.Range("H85").Value = 2404.4092
$ws.Range("I85").Value = 2480.8096
$ws.Range("J85").Value = 800
$ws.Range("K85").Value = 2480.8096
$ws.Range("L85").Value = 800
$ws.Range("M85").Value = -1232.8096
$ws.Range("N85").Value = -3296
# row 93 (hunk 33)
$ws.Range("H93").Value = 51006.375
$ws.Range("I93").Value = 1150.1428
$ws.Range("K93").Value = 1150.1428
$ws.Range("M93").Value = 97.85719999999992
# row 98 (hunk 34)
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents() | Out-Null
# row 103 (hunk 35)
$ws.Range("H103").Value = 22167
$ws.Range("J103").Value = 22167
$ws.Range("L103").Value = 22167
$ws.Range("N103").Value = -24511

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 51 (hunk 36)
$ws.Range("H51").Value = 13000
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents() | Out-Null
# row 52 (hunk 37)
$ws.Range("H52").Value = 8772.25
$ws.Range("J52").Value = 11363
$ws.Range("L52").Value = 11363
$ws.Range("N52").Value = -11815
# row 81 (hunk 38)
$ws.Range("H81").Value = 1000
$ws.Range("J81").Value = 1000
$ws.Range("L81").Value = 2000
$ws.Range("N81").Value = -4122
# row 84 (hunk 39)
$ws.Range("H84").Value = 1000
$ws.Range("J84").Value = 1000
$ws.Range("L84").Value = 10000
$ws.Range("N84").Value = -20608
# row 95 (hunk 40)
$ws.Range("H95").Value = 49999.5
$ws.Range("J95").Value = 49999.5
$ws.Range("L95").Value = 49999.5
$ws.Range("N95").Value = -55491.5
